# Edits made during key documents generation
#
# Source workbook has columns A (labels) / B (values) describing a
# MOSFET heatsink calculation. This edit:
#   1. Updates the "Power Dissipated" input (B18) from 12 to 14 W.
#   2. Duplicates the whole A:B calculation block into a new D:E block
#      for a second MOSFET ("IPP045N10N3GXKSA1") being evaluated at a
#      400lpm airflow condition, with its own (slightly different)
#      input values, re-using the same formula shapes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing Power Dissipated value -------------------------
$ws.Range("B18").Value = 14

# --- 2. Build the new D:E block -----------------------------------------
$ws.Range("D1:E1").Merge() | Out-Null
$ws.Range("D1").Value = "Heatsink for MOSFET at 400lpm"

$ws.Range("D2").Value = "MOSFET INFO"
$ws.Range("E2").Value = "IPP045N10N3GXKSA1"

$ws.Range("D3").Value = "Width"
$ws.Range("E3").Value = 0.0097

$ws.Range("D4").Value = "Length"
$ws.Range("E4").Value = 0.01484

$ws.Range("D5").Value = "Thermal Resistance, (J-MB)"
$ws.Range("E5").Value = 0.7

$ws.Range("D7").Value = "Thermal Pad Info"
$ws.Range("E7").Value = "BER220-ND"

$ws.Range("D8").Value = "Area"
$ws.Range("E8").Formula = "=E3*E4"

$ws.Range("D9").Value = "Thermal Conductivity "
$ws.Range("E9").Value = 0.9

$ws.Range("D10").Value = "Thickness"
$ws.Range("E10").Value = 0.000152

$ws.Range("D11").Value = "Thermal Resistance"
$ws.Range("E11").Formula = "=E10/(E8*E9)"

$ws.Range("D13").Value = "Heatsink Info"

$ws.Range("D14").Value = "Delta T"
$ws.Range("E14").Value = 60

$ws.Range("D15").Value = "Power"
$ws.Range("E15").Value = 12

$ws.Range("D16").Value = "Thermal Resistance"
$ws.Range("E16").Value = 2.5

$ws.Range("D18").Value = "Power Dissipated"
$ws.Range("E18").Value = 14

$ws.Range("D19").Value = "Ambient Temperature"
$ws.Range("E19").Value = 40

$ws.Range("D20").Value = "Temperature"
$ws.Range("E20").Formula = "=E19+E18*(E16+E11+E5)"

# --- 3. Formatting ---------------------------------------------------------
# Center-align the two title rows and the new part-number cell, matching
# the look of the original A1:B1 / B2 formatting.
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("E2").HorizontalAlignment = -4108

# Scientific notation for the width/length/thickness inputs, matching B3/B4/B10.
$ws.Range("E3").NumberFormat = "0.00E+00"
$ws.Range("E4").NumberFormat = "0.00E+00"
$ws.Range("E10").NumberFormat = "0.00E+00"

# Column widths: widen B to fit the new layout, and size D/E similarly to A.
$ws.Columns.Item(2).ColumnWidth = 28.666666666666668
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668
$ws.Columns.Item(5).ColumnWidth = 21.833333333333332

# --- 4. Selection / view ----------------------------------------------------
$ws.Range("D2:E20").Select() | Out-Null
